$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K column (최종점수) values by +0.1
$ws.Range("K2").Value = 59.7
$ws.Range("K3").Value = 58.5
$ws.Range("K4").Value = 50.5
$ws.Range("K5").Value = 49.5
$ws.Range("K6").Value = 46.5

# Update N column (MACRO_SCORE) values to the new recalculated constant
$ws.Range("N2").Value = 54.83846622768671
$ws.Range("N3").Value = 54.83846622768671
$ws.Range("N4").Value = 54.83846622768671
$ws.Range("N5").Value = 54.83846622768671
$ws.Range("N6").Value = 54.83846622768671
